$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C ("Förändrad") holds the "last changed" date serial for every
# data row (2 through 115). The whole column was bumped by one day
# (45180 -> 45181) in this automatic update.
$ws.Range("C2:C115").Value = 45181
